# Append two new participant rows (18 and 19) to the scores sheet,
# mirroring the existing pattern (ID, group, member, score_1, score_2,
# total_score formula), and move the active selection to D18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: Republican, member 2, score_1 3, score_2 6 -> total 9
$ws.Range("A18").Value = "6687bfc6b173fbe99a38cfb4"
$ws.Range("B18").Value = "Republican"
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 6
$ws.Range("F18").Formula = "=D18+E18"

# Row 19: Democrat, member 2, score_1 1, score_2 4 -> total 5
$ws.Range("A19").Value = "667ee4bc53bc8bcecb7843a9"
$ws.Range("B19").Value = "Democrat"
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 4
$ws.Range("F19").Formula = "=D19+E19"

# Move the selection to match the post-edit state.
$ws.Range("D18").Select()
